$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.141.18'
$ws.Range("E2").Value = '  -3.39%  '

# Row 3
$ws.Range("D3").Value = '3.513.12'
$ws.Range("E3").Value = '  -5.52%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").Value = '''579.91'
$ws.Range("E5").Value = '  -0.97%  '

# Row 6
$ws.Range("D6").Value = '''171.45'
$ws.Range("E6").Value = '  -4.90%  '

# Row 7
$ws.Range("D7").Value = '''0.609'
$ws.Range("E7").Value = '  -1.13%  '

# Row 8
$ws.Range("D8").Value = '3.509.61'
$ws.Range("E8").Value = '  -5.57%  '

# Row 9
$ws.Range("E9").Value = '  +0.15%  '

# Row 10
$ws.Range("D10").Value = '''0.189'
$ws.Range("E10").Value = '  -6.51%  '

# Row 11
$ws.Range("D11").Value = '''6.79'
$ws.Range("E11").Value = '  +4.98%  '

# Row 12
$ws.Range("D12").Value = '''0.586'
$ws.Range("E12").Value = '  -4.46%  '

# Row 13
$ws.Range("D13").Value = '''46.59'
$ws.Range("E13").Value = '  -6.13%  '

# Row 14
$ws.Range("D14").Value = '''0.0000273'
$ws.Range("E14").Value = '  -5.17%  '

# Row 15
$ws.Range("D15").Value = '4.079.06'
$ws.Range("E15").Value = '  -5.28%  '

# Row 16
$ws.Range("D16").Value = '''636.05'
$ws.Range("E16").Value = '  -6.53%  '

# Row 17
$ws.Range("E17").Value = '  -5.89%  '

# Row 18
$ws.Range("D18").Value = '69.185.67'
$ws.Range("E18").Value = '  -3.45%  '

# Row 19
$ws.Range("D19").Value = '3.512.50'
$ws.Range("E19").Value = '  -4.74%  '

# Row 20
$ws.Range("E20").Value = '  -1.02%  '

# Row 21
$ws.Range("D21").Value = '''17.36'
$ws.Range("E21").Value = '  -4.08%  '

# Row 22
$ws.Range("D22").Value = '''11.12'
$ws.Range("E22").Value = '  -4.66%  '

# Row 23
$ws.Range("D23").Value = '''0.890'
$ws.Range("E23").Value = '  -6.00%  '

# Row 24
$ws.Range("D24").Value = '''16.07'
$ws.Range("E24").Value = '  -8.30%  '

# Row 25
$ws.Range("D25").Value = '''97.68'
$ws.Range("E25").Value = '  -4.99%  '

# Row 26
$ws.Range("D26").Value = '''3.83'
$ws.Range("E26").Value = '  -4.66%  '

# Row 27
$ws.Range("E27").Value = '  -0.05%  '

# Row 28
$ws.Range("D28").Value = '''2.64'
$ws.Range("E28").Value = '  -7.26%  '

# Row 29
$ws.Range("D29").Value = '''9.34'
$ws.Range("E29").Value = '  -10.79%  '

# Row 30
$ws.Range("D30").Value = '''32.77'
$ws.Range("E30").Value = '  -8.39%  '

# Row 31
$ws.Range("D31").Value = '''3.18'
$ws.Range("E31").Value = '  -8.03%  '

# Row 32
$ws.Range("D32").Value = '''8.56'
$ws.Range("E32").Value = '  -7.26%  '

# Row 33
$ws.Range("E33").Value = '  -9.19%  '

# Row 34
$ws.Range("D34").Value = '''7.06'
$ws.Range("E34").Value = '  -4.68%  '

# Row 35
$ws.Range("D35").Value = '''619.20'
$ws.Range("E35").Value = '  +5.31%  '

# Row 36
$ws.Range("D36").Value = '''10.79'

# Row 37
$ws.Range("D37").Value = '''3.53'
$ws.Range("E37").Value = '  -14.21%  '

# Row 38
$ws.Range("D38").Value = '''0.104'
$ws.Range("E38").Value = '  -4.99%  '

# Row 39
$ws.Range("D39").Value = '''56.84'
$ws.Range("E39").Value = '  -4.06%  '

# Row 40
$ws.Range("E40").Value = '  +0.10%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '''0.0447'
$ws.Range("E41").Value = '  -2.89%  '

# Row 42
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '''0.138'
$ws.Range("E42").Value = '  -6.31%  '

# Row 43
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '3.378.50'
$ws.Range("E43").Value = '  -8.50%  '

# Row 44
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = '''0.330'
$ws.Range("E44").Value = '  -5.90%  '

# Row 45
$ws.Range("D45").Value = '''33.00'
$ws.Range("E45").Value = '  -7.75%  '

# Row 46
$ws.Range("D46").Value = '0.0₃0699'
$ws.Range("E46").Value = '  -10.75%  '

# Row 47
$ws.Range("E47").Value = '  -8.49%  '

# Row 48
$ws.Range("D48").Value = '''2.79'
$ws.Range("E48").Value = '  -3.52%  '

# Row 49
$ws.Range("E49").Value = '  -2.58%  '

# Row 50
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '''132.34'
$ws.Range("E50").Value = '  -2.92%  '

# Row 51
$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").Value = '''5.69'
$ws.Range("E51").Value = '  +15.52%  '

